# ============================================================================
# Update gh-pages "杭州-漫展信息" workbook to the data snapshot generated at
# commit 456a3b4:
#   - refresh "想去人数" (interested-count, column F) on sheet "展览" (Exhibition)
#   - refresh the same counters, plus the newly scraped/merged events for
#     early May and early June, on sheet "全部类型" (All types) so it stays in
#     sync with the per-category sheets.
# ============================================================================

$wb = $excel.ActiveWorkbook
$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll = $wb.Worksheets.Item("全部类型")

# ---- "展览": bump "想去人数" (F) counts ----
$wsExhibit.Range("F3").Value = 555
$wsExhibit.Range("F5").Value = 511
$wsExhibit.Range("F6").Value = 1148
$wsExhibit.Range("F7").Value = 330
$wsExhibit.Range("F9").Value = 123
$wsExhibit.Range("F14").Value = 867
$wsExhibit.Range("F15").Value = 863
$wsExhibit.Range("F17").Value = 67
$wsExhibit.Range("F20").Value = 757
$wsExhibit.Range("F21").Value = 1736
$wsExhibit.Range("F22").Value = 2790
$wsExhibit.Range("F23").Value = 798
$wsExhibit.Range("F24").Value = 84
$wsExhibit.Range("F25").Value = 2128
$wsExhibit.Range("F26").Value = 665
$wsExhibit.Range("F27").Value = 2975
$wsExhibit.Range("F28").Value = 571
$wsExhibit.Range("F32").Value = 723
$wsExhibit.Range("F34").Value = 128
$wsExhibit.Range("F36").Value = 1047
$wsExhibit.Range("F41").Value = 181
$wsExhibit.Range("F43").Value = 170

# ---- "全部类型": bump "想去人数" (F) counts for rows untouched by the merge below ----
$wsAll.Range("F3").Value = 555
$wsAll.Range("F5").Value = 511
$wsAll.Range("F6").Value = 1148
$wsAll.Range("F7").Value = 330
$wsAll.Range("F9").Value = 123
$wsAll.Range("F13").Value = 867
$wsAll.Range("F14").Value = 863
$wsAll.Range("F17").Value = 67
$wsAll.Range("F20").Value = 757
$wsAll.Range("F21").Value = 1736
$wsAll.Range("F22").Value = 2790
$wsAll.Range("F23").Value = 798
$wsAll.Range("F24").Value = 84
$wsAll.Range("F36").Value = 723
$wsAll.Range("F38").Value = 128
$wsAll.Range("F45").Value = 181
$wsAll.Range("F47").Value = 170

# ---- "全部类型": rows 26-28 (early May) and 40-42 (early June) shift down as a
#      newly scraped event is merged into each date-sorted block; rewrite the
#      full row contents to match. ----

# Row 26
$wsAll.Range("C26").Value = '杭州·第四届华盟动漫次元嘉年华'
$wsAll.Range("D26").Value = '创意路1号 中国智谷富春园区'
$wsAll.Range("E26").Value = '2024.05.02 10:00-05.03 17:00'
$wsAll.Range("F26").Value = 2975
$wsAll.Range("G26").Value = 58
$wsAll.Range("H26").Value = 'https://show.bilibili.com/platform/detail.html?id=82465'
$wsAll.Range("I26").Value = '//i0.hdslb.com/bfs/openplatform/202403/4XHyqi3D1709780326858.jpeg'

# Row 27
$wsAll.Range("C27").Value = '杭州·造梦探险家Porject6野蛮冲撞——第五人格ONLY'
$wsAll.Range("D27").Value = '欢西路1号 天都城酒店'
$wsAll.Range("E27").Value = '2024.05.02 10:00-05.02 22:00'
$wsAll.Range("F27").Value = 571
$wsAll.Range("G27").Value = 28
$wsAll.Range("H27").Value = 'https://show.bilibili.com/platform/detail.html?id=82851'
$wsAll.Range("I27").Value = '//i1.hdslb.com/bfs/openplatform/202403/a7IYN66u1711441126355.png'

# Row 28
$wsAll.Range("B28").NumberFormat = "@"
$wsAll.Range("B28").Value = '2024-05-03'
$wsAll.Range("C28").Value = '杭州·EY动漫嘉年华'
$wsAll.Range("D28").Value = '环丁路1428号 金色大唐城'
$wsAll.Range("E28").Value = '2024.05.03 10:00-05.05 17:00'
$wsAll.Range("F28").Value = 5
$wsAll.Range("G28").Value = 55
$wsAll.Range("H28").Value = 'https://show.bilibili.com/platform/detail.html?id=83925'
$wsAll.Range("I28").Value = '//i1.hdslb.com/bfs/openplatform/202404/0agt5uXM1712569450867.jpeg'

# Row 40
$wsAll.Range("B40").NumberFormat = "@"
$wsAll.Range("B40").Value = '2024-06-05'
$wsAll.Range("C40").Value = '杭州·英雄时代2024哈瓦西钢琴演奏会'
$wsAll.Range("D40").Value = '中国杭州北山路86号西湖岳湖景区 中国杭州西湖岳湖景区印象西湖'
$wsAll.Range("E40").Value = '2024.06.05 20:00-06.05 21:30'
$wsAll.Range("F40").Value = 1
$wsAll.Range("G40").Value = 499
$wsAll.Range("H40").Value = 'https://show.bilibili.com/platform/detail.html?id=83902'
$wsAll.Range("I40").Value = '//i2.hdslb.com/bfs/openplatform/202404/BFRFmKpT1712569969076.jpeg'

# Row 41
$wsAll.Range("B41").NumberFormat = "@"
$wsAll.Range("B41").Value = '2024-06-08'
$wsAll.Range("C41").Value = '杭州·第八届YH樱花动漫游戏文化节'
$wsAll.Range("D41").Value = '德胜东路2539号 梦马汽车小镇'
$wsAll.Range("E41").Value = '2024.06.08 10:00-06.10 17:00'
$wsAll.Range("F41").Value = 1047
$wsAll.Range("G41").Value = 65
$wsAll.Range("H41").Value = 'https://show.bilibili.com/platform/detail.html?id=82687'
$wsAll.Range("I41").Value = '//i2.hdslb.com/bfs/openplatform/202403/S5pnadXj1710210939138.png'

# Row 42
$wsAll.Range("B42").NumberFormat = "@"
$wsAll.Range("B42").Value = '2024-06-09'
$wsAll.Range("C42").Value = '杭州·第三届日夜国乙only'
$wsAll.Range("D42").Value = '创意路1号 中国智谷富春园区'
$wsAll.Range("E42").Value = '2024.06.09 10:00-06.09 23:00'
$wsAll.Range("F42").Value = 1757
$wsAll.Range("G42").Value = 58
$wsAll.Range("H42").Value = 'https://show.bilibili.com/platform/detail.html?id=82618'
$wsAll.Range("I42").Value = '//i2.hdslb.com/bfs/openplatform/202403/fXRzYEFH1710124366279.png'

Write-Output "applied 漫展信息 refresh (456a3b4)"
